$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # column G = 7
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
